# Applies the 26-11-2023 20:30 script update to the Croatia Prva NL 2023-2024
# worksheet: a handful of rows whose match details (columns F:V) were
# re-shuffled between rows sharing the same match date (A:E stay put), plus
# four brand-new fixtures appended at the end of the table (rows 93-96).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: capture the CURRENT (pre-edit) F:V contents of every row that is
# part of a reshuffle, keyed by its current row number. All reads happen
# before any writes so a row's data is never read back after being
# overwritten.
# ---------------------------------------------------------------------------
$sourceRows = @(3,4,5,6,15,16,17,21,22,27,28,29,35,36,46,48,57,58,59,82,83,90,91)
$captured = @{}
foreach ($r in $sourceRows) {
    $captured[$r] = $ws.Range("F$r`:V$r").Value2
}

# ---------------------------------------------------------------------------
# Step 2: write each captured block back into its new row. Map: new row ->
# row whose old content now lives there.
# ---------------------------------------------------------------------------
$mapping = @{
    3  = 6
    4  = 3
    5  = 4
    6  = 5
    15 = 17
    16 = 15
    17 = 16
    21 = 22
    22 = 21
    27 = 29
    28 = 27
    29 = 28
    35 = 36
    36 = 35
    46 = 48
    48 = 46
    57 = 58
    58 = 59
    59 = 57
    82 = 83
    83 = 82
    90 = 91
    91 = 90
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $ws.Range("F$newRow`:V$newRow").Value2 = $captured[$oldRow]
}

# ---------------------------------------------------------------------------
# Step 3: append the four new fixtures as rows 93-96. First clone the
# formatting of the last existing data row (92) so styles (bold/bordered
# index column, date-formatted match-date column, etc.) match the rest of
# the table, then populate the values.
# ---------------------------------------------------------------------------
$ws.Range("A92:V92").Copy()
$ws.Range("A93:V96").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$newRows = @(
    @{ row=93; A=92; E=45255.53472222222; F="Solin";       G=0; H="Cibalia";          I=3;
       J=1.72; K="24/11/2023 01:12"; L=1.7;  M="25/11/2023 12:36";
       N=3.54; O="24/11/2023 01:12"; P=3.75; Q="25/11/2023 12:36";
       R=3.97; S="24/11/2023 01:12"; T=4.59; U="25/11/2023 12:18";
       V="https://www.betexplorer.com/football/croatia/prva-nl/solin-cibalia/EXwhrDGf/" },

    @{ row=94; A=93; E=45255.5625;         F="Bijelo Brdo"; G=0; H="Zrinski Jurjevac"; I=1;
       J=3.15; K="24/11/2023 01:42"; L=4.97; M="25/11/2023 13:26";
       N=3;    O="24/11/2023 01:42"; P=3.14; Q="25/11/2023 13:29";
       R=2.16; S="24/11/2023 01:42"; T=1.81; U="25/11/2023 13:29";
       V="https://www.betexplorer.com/football/croatia/prva-nl/bijelo-brdo-zrinski-jurjevac/67zppZor/" },

    @{ row=95; A=94; E=45255.5625;         F="Dubrava";     G=2; H="Croatia Zmijavci";  I=1;
       J=1.95; K="24/11/2023 01:42"; L=2.21; M="25/11/2023 13:17";
       N=3.38; O="24/11/2023 01:42"; P=3.1;  Q="25/11/2023 13:17";
       R=3.25; S="24/11/2023 01:42"; T=3.36; U="25/11/2023 13:17";
       V="https://www.betexplorer.com/football/croatia/prva-nl/dubrava-zagreb-croatia-zmijavci/2gKn4fwE/" },

    @{ row=96; A=95; E=45255.5625;         F="Jarun";       G=0; H="Dugopolje";         I=2;
       J=2.22; K="24/11/2023 01:42"; L=1.87; M="25/11/2023 13:26";
       N=3.21; O="24/11/2023 01:42"; P=3.67; Q="25/11/2023 13:26";
       R=2.83; S="24/11/2023 01:42"; T=3.8;  U="25/11/2023 13:26";
       V="https://www.betexplorer.com/football/croatia/prva-nl/jarun-dugopolje/zmS0tio7/" }
)

foreach ($nr in $newRows) {
    $r = $nr.row
    $ws.Range("A$r").Value2 = $nr.A
    $ws.Range("B$r").Value2 = "croatia"
    $ws.Range("C$r").Value2 = "prva-nl"
    $ws.Range("D$r").Value2 = "2023-2024"
    $ws.Range("E$r").Value2 = $nr.E
    $ws.Range("E$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("F$r").Value2 = $nr.F
    $ws.Range("G$r").Value2 = $nr.G
    $ws.Range("H$r").Value2 = $nr.H
    $ws.Range("I$r").Value2 = $nr.I
    $ws.Range("J$r").Value2 = $nr.J
    $ws.Range("K$r").Value2 = $nr.K
    $ws.Range("L$r").Value2 = $nr.L
    $ws.Range("M$r").Value2 = $nr.M
    $ws.Range("N$r").Value2 = $nr.N
    $ws.Range("O$r").Value2 = $nr.O
    $ws.Range("P$r").Value2 = $nr.P
    $ws.Range("Q$r").Value2 = $nr.Q
    $ws.Range("R$r").Value2 = $nr.R
    $ws.Range("S$r").Value2 = $nr.S
    $ws.Range("T$r").Value2 = $nr.T
    $ws.Range("U$r").Value2 = $nr.U
    $ws.Range("V$r").Value2 = $nr.V
}

# ---------------------------------------------------------------------------
# Step 4: make sure the sheet's used-range dimension now reaches row 96
# (setting a value outside the previous dimension already grows it, but
# touch the corner cell explicitly so it's unambiguous).
# ---------------------------------------------------------------------------
$null = $ws.Range("A1:V96")
